# Auto-generated: apply 2025-11-18 daily crime-count update to violent-crime-full-year workbook
# Source data: citywide totals + per-neighborhood sheets + the "By Neighborhood" rollup sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5899
$ws.Range("L3").Value = 6419
$ws.Range("B4").Value = 1719
$ws.Range("K4").Value = 1795
$ws.Range("L4").Value = 1582
$ws.Range("L6").Value = 5276
$ws.Range("B7").Value = 23351
$ws.Range("K7").Value = 27587
$ws.Range("L7").Value = 19559

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 169
$ws.Range("L5").Value = 70
$ws.Range("L6").Value = 151
$ws.Range("L7").Value = 631
$ws.Range("L8").Value = 1286
$ws.Range("L10").Value = 130
$ws.Range("L11").Value = 325
$ws.Range("L15").Value = 159
$ws.Range("L20").Value = 496
$ws.Range("L23").Value = 212
$ws.Range("K25").Value = 131
$ws.Range("L29").Value = 1096
$ws.Range("L33").Value = 885
$ws.Range("L36").Value = 246
$ws.Range("L42").Value = 625
$ws.Range("L43").Value = 145
$ws.Range("L49").Value = 104
$ws.Range("L51").Value = 248
$ws.Range("L53").Value = 215
$ws.Range("L54").Value = 425
$ws.Range("L60").Value = 127
$ws.Range("B63").Value = 423
$ws.Range("L63").Value = 60
$ws.Range("L67").Value = 674
$ws.Range("L76").Value = 298
$ws.Range("L77").Value = 132
$ws.Range("L78").Value = 250
$ws.Range("L79").Value = 545
$ws.Range("L80").Value = 66
$ws.Range("L83").Value = 427
$ws.Range("L85").Value = 969
$ws.Range("L88").Value = 208
$ws.Range("L91").Value = 263
$ws.Range("L92").Value = 60
$ws.Range("L94").Value = 242
$ws.Range("L95").Value = 275
$ws.Range("L97").Value = 160
$ws.Range("L99").Value = 341
$ws.Range("L100").Value = 36
$ws.Range("B101").Value = 23351
$ws.Range("K101").Value = 27587
$ws.Range("L101").Value = 19559

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 214
$ws.Range("L3").Value = 201
$ws.Range("L4").Value = 46
$ws.Range("L7").Value = 631

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 122
$ws.Range("L7").Value = 325

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 292
$ws.Range("L3").Value = 399
$ws.Range("L6").Value = 202
$ws.Range("L7").Value = 969

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 457
$ws.Range("L6").Value = 315
$ws.Range("L7").Value = 1286

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 136
$ws.Range("L7").Value = 427

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 310
$ws.Range("L4").Value = 59
$ws.Range("L6").Value = 255
$ws.Range("L7").Value = 885

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 102
$ws.Range("L7").Value = 275

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 74
$ws.Range("L7").Value = 341

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 263
$ws.Range("L7").Value = 674

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 104
$ws.Range("L7").Value = 425

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L4").Value = 59
$ws.Range("L7").Value = 1096

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 63
$ws.Range("L7").Value = 298

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 170
$ws.Range("L7").Value = 625

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 66
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 212

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 88
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 147
$ws.Range("L7").Value = 545

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 155
$ws.Range("L3").Value = 172
$ws.Range("L7").Value = 496

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 79
$ws.Range("L7").Value = 246

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 242

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 60
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 159

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 56
$ws.Range("L3").Value = 54
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L3").Value = 71
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 77
$ws.Range("L7").Value = 248

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 66
